# Fix feared future self formulation
#
# The "feared future self" activity wording is reworded from
# "your/my feared future self if you ..." to
# "the feared future self you/I might become if you/I ..."
# for both the Smoking ("continue to smoke") and the Physical Activity
# ("fail to become more physically active") variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activities")

# --- Physical activity rows (18th and 19th activity: "fail to become more physically active") ---
$ws.Range("E19").Value = "thinking about the feared future self you might become if you fail to become more physically active"
$ws.Range("F19").Value = "Thinking about the feared future self you might become if you fail to become more physically active"
$ws.Range("G19").Value = "think about the feared future self you might become if you fail to become more physically active"
$ws.Range("H19").Value = "think about the feared future self I might become if I fail to become more physically active"

$ws.Range("E20").Value = "thinking about the feared future self you might become if you fail to become more physically active"
$ws.Range("F20").Value = "Thinking about the feared future self you might become if you fail to become more physically active"
$ws.Range("G20").Value = "think about the feared future self you might become if you fail to become more physically active"
$ws.Range("H20").Value = "think about the feared future self I might become if I fail to become more physically active"

# --- Smoking rows (3rd and 4th activity: "continue to smoke") ---
$ws.Range("F4").Value = "Thinking about the feared future self you might become if you continue to smoke"
$ws.Range("E4").Value = "thinking about the feared future self you might become if you continue to smoke"
$ws.Range("G4").Value = "think about the feared future self you might become if you continue to smoke"
$ws.Range("H4").Value = "think about the feared future self I might become if I continue to smoke"

$ws.Range("F5").Value = "Thinking about the feared future self you might become if you continue to smoke"
$ws.Range("E5").Value = "thinking about the feared future self you might become if you continue to smoke"
$ws.Range("G5").Value = "think about the feared future self you might become if you continue to smoke"
$ws.Range("H5").Value = "think about the feared future self I might become if I continue to smoke"

# --- Leave the selection where the author ended up after editing row 19 ---
$ws.Range("H19").Select()
